$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "64.437.33"
$ws.Range("E2").Value = "  +0.93%  "

Set-TextValue $ws.Range("D3") "3.111.39"
$ws.Range("E3").Value = "  -1.60%  "

Set-TextValue $ws.Range("D4") "1.01"
$ws.Range("E4").Value = "  +0.37%  "

Set-TextValue $ws.Range("D5") "589.53"
$ws.Range("E5").Value = "  +0.67%  "

Set-TextValue $ws.Range("D6") "152.02"
$ws.Range("E6").Value = "  +4.33%  "

Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.04%  "

Set-TextValue $ws.Range("D8") "3.108.30"
$ws.Range("E8").Value = "  -0.78%  "

Set-TextValue $ws.Range("D9") "0.530"
$ws.Range("E9").Value = "  +0.52%  "

Set-TextValue $ws.Range("D10") "0.159"
$ws.Range("E10").Value = "  -0.39%  "

Set-TextValue $ws.Range("D11") "5.94"
$ws.Range("E11").Value = "  +2.88%  "

Set-TextValue $ws.Range("D12") "0.460"
$ws.Range("E12").Value = "  +0.92%  "

Set-TextValue $ws.Range("D13") "37.84"
$ws.Range("E13").Value = "  +2.92%  "

Set-TextValue $ws.Range("D14") "0.0000243"
$ws.Range("E14").Value = "  -1.14%  "

Set-TextValue $ws.Range("D15") "3.628.83"
$ws.Range("E15").Value = "  -1.49%  "

$ws.Range("E16").Value = "  -1.64%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D17") "7.22"
$ws.Range("E17").Value = "  +2.58%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D18") "63.969.25"
$ws.Range("E18").Value = "  +0.49%  "

Set-TextValue $ws.Range("D19") "3.108.45"
$ws.Range("E19").Value = "  -1.26%  "

Set-TextValue $ws.Range("D20") "468.48"
$ws.Range("E20").Value = "  +1.32%  "

Set-TextValue $ws.Range("D21") "14.89"
$ws.Range("E21").Value = "  +4.28%  "

Set-TextValue $ws.Range("D22") "0.739"
$ws.Range("E22").Value = "  +1.06%  "

Set-TextValue $ws.Range("D23") "7.58"
$ws.Range("E23").Value = "  +2.54%  "

Set-TextValue $ws.Range("D24") "13.27"
$ws.Range("E24").Value = "  +3.02%  "

Set-TextValue $ws.Range("D25") "2.37"
$ws.Range("E25").Value = "  +7.11%  "

Set-TextValue $ws.Range("D26") "81.68"
$ws.Range("E26").Value = "  +1.12%  "

Set-TextValue $ws.Range("D27") "0.999"
$ws.Range("E27").Value = "  -0.18%  "

Set-TextValue $ws.Range("D28") "9.83"
$ws.Range("E28").Value = "  +6.38%  "

$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D29") "7.42"
$ws.Range("E29").Value = "  +4.54%  "

Set-TextValue $ws.Range("D30") "2.69"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D31") "1.01"
$ws.Range("E31").Value = "  +0.49%  "

Set-TextValue $ws.Range("D32") "2.21"
$ws.Range("E32").Value = "  +0.38%  "

Set-TextValue $ws.Range("D33") "0.117"
$ws.Range("E33").Value = "  +7.53%  "

Set-TextValue $ws.Range("D34") "27.41"
$ws.Range("E34").Value = "  +2.30%  "

Set-TextValue $ws.Range("D35") "0.0₃0851"
$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("E36").Value = "  +1.46%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D37") "3.38"
$ws.Range("E37").Value = "  +1.43%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D38") "6.15"
$ws.Range("E38").Value = "  +2.78%  "

Set-TextValue $ws.Range("D39") "2.27"
$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D40") "9.37"
$ws.Range("E40").Value = "  +5.92%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D41") "50.81"
$ws.Range("E41").Value = "  -0.99%  "

Set-TextValue $ws.Range("D42") "452.65"
$ws.Range("E42").Value = "  +4.14%  "

Set-TextValue $ws.Range("D43") "0.292"
$ws.Range("E43").Value = "  +4.04%  "

Set-TextValue $ws.Range("D44") "0.0371"
$ws.Range("E44").Value = "  +0.23%  "

Set-TextValue $ws.Range("D45") "2.850.89"
$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("E46").Value = "  +1.96%  "

Set-TextValue $ws.Range("D47") "38.76"
$ws.Range("E47").Value = "  +3.48%  "

Set-TextValue $ws.Range("D48") "130.03"
$ws.Range("E48").Value = "  +2.67%  "

$ws.Range("E49").Value = "  +0.04%  "

Set-TextValue $ws.Range("D50") "25.19"
$ws.Range("E50").Value = "  +4.98%  "

Set-TextValue $ws.Range("D51") "2.27"
$ws.Range("E51").Value = "  +4.61%  "
